$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp 1 - Classifiers test")

# New header labels (B1 changes from "Accuracy" to "Accuracy in CV"; new C1 "Final accuracy")
$ws.Range("B1").Value = "Accuracy in CV"
$ws.Range("C1").Value = "Final accuracy"

# New column C width
$ws.Columns.Item(3).ColumnWidth = 14.2

# Make this sheet the active / selected tab, with C2 selected
$ws.Activate()
$ws.Range("C2").Select()
